# feat: add 2022-Q3 data
#
# Before: 总计, 2022-Q2, 2021-Q4
# After : 总计, 2022-Q3 (new), 2022-Q2, 2021-Q4
#
# The new "2022-Q3" tab carries the fund-holding snapshot that used to
# live on the "2022-Q2" tab, refreshed with the new quarter's numbers.
# The old "2022-Q2" tab keeps its original numbers untouched, and
# "2021-Q4" is untouched as well. The "总计" (summary) sheet gets a new
# row for 2022-Q3.

$wb = $excel.ActiveWorkbook

# 1) Insert the brand-new "2022-Q3" sheet right before the existing
#    "2022-Q2" sheet (i.e. right after "总计"), then rename it.
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$wb.Worksheets.Add($beforeSheet)
$wb.Worksheets.Item(2).Name = "2022-Q3"

# NOTE: always look sheets back up by name/index AFTER Worksheets.Add()
# runs -- handles captured before an Add() call don't reliably receive
# .Copy() pastes anymore.
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# 2) Seed the new sheet from the old "2022-Q2" layout (header row + the
#    single fund row) so styles/borders match, then overwrite with this
#    quarter's refreshed figures.
$q2Sheet.Range("A1:H2").Copy($q3Sheet.Range("A1"))

$q3Sheet.Range("D2").Value = "'6.25"
$q3Sheet.Range("E2").Value = "'62.12"
$q3Sheet.Range("F2").Value = "'4.41"
$q3Sheet.Range("G2").Value = "'0.2756"
$q3Sheet.Range("H2").Value = 3

# 3) Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    above the existing 2022-Q2 row, pushing 2022-Q2/2021-Q4 down.
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.28

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

Write-Output "done"
